# Commit: "changed value of cells in sample file"
#
# The sample data-collection workbook had every manually-entered employee
# count bumped from 1 to 2 on the four data-entry tables. The "Totals Check"
# and "Matching Check" sheets reference these cells via formulas, so they
# recompute automatically once the inputs change.

$wb = $excel.ActiveWorkbook

# 1. Number of Employees -- data rows are 7-16, columns B:Q
$wsEmployees = $wb.Worksheets.Item("1.Number of Employees")
$wsEmployees.Range("B7:Q16").Value = 2

# 2. Compensation -- data rows are 6-15, columns B:Q
$wsComp = $wb.Worksheets.Item("2.Compensation")
$wsComp.Range("B6:Q15").Value = 2

# 3. Performance Pay -- data rows are 6-15, columns B:Q
$wsPerf = $wb.Worksheets.Item("3.Performance Pay")
$wsPerf.Range("B6:Q15").Value = 2

# 4. Tenure -- data rows are 6-15, columns B:Q
$wsTenure = $wb.Worksheets.Item("4.Tenure")
$wsTenure.Range("B6:Q15").Value = 2

# The workbook was left with "1.Number of Employees" as the active tab and
# A26 selected there (previously "2.Compensation" / B6:Q15 was active).
$wsEmployees.Activate()
$wsEmployees.Range("A26").Select()

Write-Output "updated cell values on the four data-entry tables"
